$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 468.83334
$ws.Range("I53").Value = 262.10526
$ws.Range("K53").Value = 262.10526
$ws.Range("M53").Value = 374.89474
# Row 86
$ws.Range("H86").Value = 5999
$ws.Range("I86").Value = 1999
$ws.Range("K86").Value = 1999
$ws.Range("M86").Value = -876
# Row 89
$ws.Range("H89").Value = 5999
$ws.Range("I89").Value = 1999
$ws.Range("K89").Value = 9995
$ws.Range("M89").Value = -4379
# Row 92
$ws.Range("H92").Value = 2048.721
$ws.Range("I92").Value = 1875.2778
$ws.Range("J92").Value = 2940.7144
$ws.Range("K92").Value = 1875.2778
$ws.Range("L92").Value = 2940.7144
$ws.Range("M92").Value = -627.2778000000001
$ws.Range("N92").Value = -5436.7144
# Row 98
$ws.Range("H98").Value = 4209114
$ws.Range("I98").Value = 4548437.5
$ws.Range("K98").Value = 4548437.5
$ws.Range("M98").Value = -4546939.5
# Row 122
$ws.Range("H122").Value = 4209114
$ws.Range("I122").Value = 4548437.5
$ws.Range("K122").Value = 13645312.5
$ws.Range("M122").Value = -13642862.5
# Row 132
$ws.Range("H132").Value = 2935.658
$ws.Range("I132").Value = 1163.1177
$ws.Range("K132").Value = 3489.3531
$ws.Range("M132").Value = -959.3531000000003
# Row 137
$ws.Range("H137").Value = 72045.125
$ws.Range("I137").Value = 140263.25
$ws.Range("J137").Value = 3827
$ws.Range("K137").Value = 420789.75
$ws.Range("L137").Value = 11481
$ws.Range("M137").Value = -418239.75
$ws.Range("N137").Value = -16581

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
# Row 45
$ws.Range("H45").Value = 1632.7142
$ws.Range("J45").Value = 1702.25
$ws.Range("L45").Value = 1702.25
$ws.Range("N45").Value = -2456.25
# Row 74
$ws.Range("H74").Value = 88564.30499999999
$ws.Range("I74").Value = 118025.88
$ws.Range("J74").Value = 5089.8335
$ws.Range("K74").Value = 118025.88
$ws.Range("L74").Value = 5089.8335
$ws.Range("M74").Value = -117151.88
$ws.Range("N74").Value = -6837.8335
# Row 77
$ws.Range("H77").Value = 88564.30499999999
$ws.Range("I77").Value = 118025.88
$ws.Range("J77").Value = 5089.8335
$ws.Range("K77").Value = 590129.4
$ws.Range("L77").Value = 25449.1675
$ws.Range("M77").Value = -585761.4
$ws.Range("N77").Value = -34185.1675
# Row 102
$ws.Range("H102").Value = 1876.5333
$ws.Range("I102").Value = 1962.6666
$ws.Range("J102").Value = 1532
$ws.Range("K102").Value = 1962.6666
$ws.Range("L102").Value = 1532
$ws.Range("M102").Value = -340.6666
$ws.Range("N102").Value = -4776
# Row 110
$ws.Range("H110").Value = 8700.700000000001
$ws.Range("I110").Value = 1401.8
$ws.Range("K110").Value = 1401.8
$ws.Range("M110").Value = 643.2
# Row 113
$ws.Range("H113").Value = 79998
$ws.Range("J113").Value = 79998
$ws.Range("L113").Value = 79998
$ws.Range("N113").Value = -88676
# Row 122
$ws.Range("H122").Value = 3079.2307
$ws.Range("I122").Value = 3064.6365
$ws.Range("K122").Value = 9193.9095
$ws.Range("M122").Value = -6743.9095

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 5683.3
$ws.Range("I20").Value = 3984.8
$ws.Range("K20").Value = 3984.8
$ws.Range("M20").Value = -3737.8
# Row 86
$ws.Range("H86").Value = 1932.2
$ws.Range("I86").Value = 1956.4706
$ws.Range("J86").Value = 1794.6666
$ws.Range("K86").Value = 1956.4706
$ws.Range("L86").Value = 1794.6666
$ws.Range("M86").Value = -833.4706000000001
$ws.Range("N86").Value = -4040.6666
# Row 89
$ws.Range("H89").Value = 1932.2
$ws.Range("I89").Value = 1956.4706
$ws.Range("J89").Value = 1794.6666
$ws.Range("K89").Value = 9782.353000000001
$ws.Range("L89").Value = 8973.333000000001
$ws.Range("M89").Value = -4166.353000000001
$ws.Range("N89").Value = -20205.333
# Row 94
$ws.Range("H94").Value = 970.3125
$ws.Range("I94").Value = 984
$ws.Range("K94").Value = 984
$ws.Range("M94").Value = -533
# Row 99
$ws.Range("H99").Value = 4878.0625
$ws.Range("I99").Value = 5091.846
$ws.Range("K99").Value = 5091.846
$ws.Range("M99").Value = -3593.846
# Row 105
$ws.Range("H105").Value = 2373.8572
$ws.Range("I105").Value = 2283.6
$ws.Range("J105").Value = 2599.5
$ws.Range("K105").Value = 2283.6
$ws.Range("L105").Value = 2599.5
$ws.Range("M105").Value = -536.5999999999999
$ws.Range("N105").Value = -6093.5
# Row 107
$ws.Range("H107").Value = 1581.3334
$ws.Range("I107").Value = 1581.3334
$ws.Range("K107").Value = 1581.3334
$ws.Range("M107").Value = 338.6666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 528948.5600000001
$ws.Range("I31").Value = 669320.9399999999
$ws.Range("K31").Value = 669320.9399999999
$ws.Range("M31").Value = -669025.9399999999
# Row 34
$ws.Range("H34").Value = 528948.5600000001
$ws.Range("I34").Value = 669320.9399999999
$ws.Range("K34").Value = 669320.9399999999
$ws.Range("M34").Value = -669118.9399999999
# Row 99
$ws.Range("H99").Value = 607334.4
$ws.Range("J99").Value = 22668.715
$ws.Range("L99").Value = 22668.715
$ws.Range("N99").Value = -25664.715
# Row 126
$ws.Range("H126").Value = 607334.4
$ws.Range("J126").Value = 22668.715
$ws.Range("L126").Value = 68006.145
$ws.Range("N126").Value = -72946.145
# Row 134
$ws.Range("H134").Value = 6471.7393
$ws.Range("I134").Value = 7150.647
$ws.Range("K134").Value = 21451.941
$ws.Range("M134").Value = -18916.941

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 3855543.8
$ws.Range("I4").Value = 4575146.5
$ws.Range("K4").Value = 13725439.5
$ws.Range("M4").Value = -13725327.5
# Row 12
$ws.Range("H12").Value = 3435.524
$ws.Range("J12").Value = 5127.7856
$ws.Range("L12").Value = 15383.3568
$ws.Range("N12").Value = -15729.3568
# Row 75
$ws.Range("H75").Value = 800
$ws.Range("J75").Value = 800
$ws.Range("L75").Value = 2400
$ws.Range("N75").Value = -4396
# Row 78
$ws.Range("H78").Value = 800
$ws.Range("J78").Value = 800
$ws.Range("L78").Value = 7200
$ws.Range("N78").Value = -17184
# Row 80
$ws.Range("H80").Value = 6861.6665
$ws.Range("I80").Value = 5449.5
$ws.Range("J80").Value = 7265.143
$ws.Range("K80").Value = 16348.5
$ws.Range("L80").Value = 21795.429
$ws.Range("M80").Value = -15412.5
$ws.Range("N80").Value = -23667.429
# Row 83
$ws.Range("H83").Value = 6861.6665
$ws.Range("I83").Value = 5449.5
$ws.Range("J83").Value = 7265.143
$ws.Range("K83").Value = 49045.5
$ws.Range("L83").Value = 65386.287
$ws.Range("M83").Value = -44365.5
$ws.Range("N83").Value = -74746.287
# Row 92
$ws.Range("H92").Value = 100
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
# Row 137
$ws.Range("H137").Value = 30304812
$ws.Range("I137").Value = 1959.6
$ws.Range("J137").Value = 333333340
$ws.Range("K137").Value = 5878.799999999999
$ws.Range("L137").Value = 1000000020
$ws.Range("M137").Value = -778.7999999999993
$ws.Range("N137").Value = -1000010220
# Row 140
$ws.Range("H140").Value = 6104001.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 20
$ws.Range("H20").Value = 27499.5
$ws.Range("J20").Value = 27499.5
$ws.Range("L20").Value = 27499.5
$ws.Range("N20").Value = -27989.5
# Row 24
$ws.Range("H24").Value = 14666
$ws.Range("J24").Value = 14666
$ws.Range("L24").Value = 14666
$ws.Range("N24").Value = -15012
# Row 47
$ws.Range("H47").Value = 16677
$ws.Range("I47").Value = 20001
$ws.Range("J47").Value = 15015
$ws.Range("K47").Value = 20001
$ws.Range("L47").Value = 15015
$ws.Range("M47").Value = -19433
$ws.Range("N47").Value = -16151
# Row 80
$ws.Range("H80").Value = 5146.0625
$ws.Range("I80").Value = 4995
$ws.Range("J80").Value = 5156.1333
$ws.Range("K80").Value = 4995
$ws.Range("L80").Value = 5156.1333
$ws.Range("M80").Value = -3997
$ws.Range("N80").Value = -7152.1333
# Row 83
$ws.Range("H83").Value = 5146.0625
$ws.Range("I83").Value = 4995
$ws.Range("J83").Value = 5156.1333
$ws.Range("K83").Value = 24975
$ws.Range("L83").Value = 25780.6665
$ws.Range("M83").Value = -19983
$ws.Range("N83").Value = -35764.66650000001
# Row 113
$ws.Range("H113").Value = 12390.714
$ws.Range("I113").Value = 1567.6
$ws.Range("K113").Value = 1567.6
$ws.Range("M113").Value = 602.4000000000001
# Row 122
$ws.Range("H122").Value = 2900.25
$ws.Range("I122").Value = 3418.6924
$ws.Range("J122").Value = 2450.9333
$ws.Range("K122").Value = 10256.0772
$ws.Range("L122").Value = 7352.7999
$ws.Range("M122").Value = -7806.0772
$ws.Range("N122").Value = -12252.7999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 4999
$ws.Range("I22").Value = 4998
$ws.Range("K22").Value = 4998
$ws.Range("M22").Value = -4703
# Row 27
$ws.Range("H27").Value = 4999
$ws.Range("I27").Value = 4998
$ws.Range("K27").Value = 4998
$ws.Range("M27").Value = -4891

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 431824
$ws.Range("I2").Value = 431824
$ws.Range("K2").Value = 431824
$ws.Range("M2").Value = -431712
# Row 4
$ws.Range("H4").Value = 4262
$ws.Range("I4").Value = 3999.3333
$ws.Range("K4").Value = 3999.3333
$ws.Range("M4").Value = -3886.3333
# Row 17
$ws.Range("H17").Value = 25370
$ws.Range("I17").Value = 25370
$ws.Range("K17").Value = 25370
$ws.Range("M17").Value = -25198
# Row 31
$ws.Range("H31").Value = 9999
$ws.Range("J31").Value = 9999
$ws.Range("L31").Value = 9999
$ws.Range("N31").Value = -10695
# Row 122
$ws.Range("H122").Value = 2398.2
$ws.Range("I122").Value = 2398.2
$ws.Range("K122").Value = 7194.599999999999
$ws.Range("M122").Value = -4744.599999999999
